# Daily attendance processing - 2025-11-17 15:23:48
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Row 2 - reorder "Recorded By" list
$ws.Range("G2").Value = "System, servinaz@med.asu.edu.eg, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 3 - reorder "Recorded By" list
$ws.Range("G3").Value = "System, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"

# Row 4 - add recorders and update attendance count
$ws.Range("G4").Value = "majorelle.magdy@med.asu.edu.eg, servinaz@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("H4").Value = "48/251"

# Updated average attendance percentages (stored as plain text, not a
# percentage number, so force text entry then restore General display)
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "24.3%"
$ws.Range("L10").NumberFormat = "General"

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "24.3%"
$ws.Range("S15").NumberFormat = "General"

# Row 28 - reorder "Recorded By" list
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
